$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D hold price text that Excel COM would otherwise auto-coerce to
# a floating-point number (e.g. "320.08" -> 320.07999999999998) and stamp a
# new number-format style on the cell. Force text entry via a temporary "@"
# (Text) number format, then restore the cell to the default "Normal" style
# so no spurious style index is left attached to the cell.
function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "48.061.24"
$ws.Range("E2").Value = "  -0.42%  "
Set-TextValue $ws.Range("D3") "2.496.39"
$ws.Range("E3").Value = "  -0.88%  "
$ws.Range("E4").Value = "  -0.04%  "
Set-TextValue $ws.Range("D5") "320.08"
$ws.Range("E5").Value = "  -1.03%  "
Set-TextValue $ws.Range("D6") "105.69"
$ws.Range("E6").Value = "  -3.26%  "
Set-TextValue $ws.Range("D7") "0.521"
$ws.Range("E7").Value = "  -0.97%  "
$ws.Range("E8").Value = "  -0.06%  "
$ws.Range("E9").Value = "  -4.29%  "
$ws.Range("E10").Value = "  -3.87%  "
Set-TextValue $ws.Range("D11") "20.05"
$ws.Range("E11").Value = "  +0.93%  "
$ws.Range("E12").Value = "  -1.98%  "
$ws.Range("E13").Value = "  -0.43%  "
$ws.Range("E14").Value = "  -2.16%  "
Set-TextValue $ws.Range("D15") "2.887.66"
$ws.Range("E15").Value = "  -0.88%  "
Set-TextValue $ws.Range("D16") "2.496.33"
$ws.Range("E16").Value = "  -1.10%  "
$ws.Range("E17").Value = "  -2.90%  "
Set-TextValue $ws.Range("D18") "47.894.51"
$ws.Range("E18").Value = "  -0.54%  "
Set-TextValue $ws.Range("D19") "13.02"
$ws.Range("E19").Value = "  -2.53%  "
Set-TextValue $ws.Range("D20") "2.96"
$ws.Range("E20").Value = "  +8.49%  "
$ws.Range("E21").Value = "  +0.00%  "
$ws.Range("E22").Value = "  -1.22%  "
Set-TextValue $ws.Range("D23") "71.05"
$ws.Range("E23").Value = "  -2.19%  "
Set-TextValue $ws.Range("D24") "271.65"
$ws.Range("E24").Value = "  +0.72%  "
$ws.Range("E25").Value = "  -3.13%  "
$ws.Range("E26").Value = "  -0.07%  "
Set-TextValue $ws.Range("D27") "25.77"
$ws.Range("E27").Value = "  -1.54%  "
Set-TextValue $ws.Range("D28") "2.28"
$ws.Range("E28").Value = "  +3.30%  "
$ws.Range("E29").Value = "  -4.53%  "
Set-TextValue $ws.Range("D30") "0.141"
$ws.Range("E30").Value = "  -3.15%  "
Set-TextValue $ws.Range("D31") "34.67"
$ws.Range("E31").Value = "  -1.52%  "
Set-TextValue $ws.Range("D32") "49.17"
$ws.Range("E32").Value = "  -1.33%  "
$ws.Range("E33").Value = "  -0.04%  "
Set-TextValue $ws.Range("D34") "19.08"
$ws.Range("E34").Value = "  -4.58%  "
$ws.Range("E35").Value = "  -2.22%  "
$ws.Range("E36").Value = "  -2.19%  "
$ws.Range("E37").Value = "  -2.58%  "
$ws.Range("E38").Value = "  -3.07%  "
$ws.Range("E39").Value = "  -4.70%  "
Set-TextValue $ws.Range("D40") "121.77"
$ws.Range("E40").Value = "  +2.82%  "
$ws.Range("E41").Value = "  -2.12%  "
Set-TextValue $ws.Range("D42") "22.13"
$ws.Range("E42").Value = "  -0.44%  "
$ws.Range("E43").Value = "  +0.87%  "
$ws.Range("E44").Value = "  +1.14%  "
Set-TextValue $ws.Range("D45") "1.998.99"
$ws.Range("E45").Value = "  -0.07%  "
$ws.Range("E46").Value = "  +0.60%  "
Set-TextValue $ws.Range("D47") "1.89"
$ws.Range("E47").Value = "  -0.51%  "
$ws.Range("E48").Value = "  -1.17%  "
Set-TextValue $ws.Range("D49") "8.90"
$ws.Range("E49").Value = "  -2.12%  "
Set-TextValue $ws.Range("D50") "5.17"
$ws.Range("E50").Value = "  -1.88%  "
Set-TextValue $ws.Range("D51") "78.75"
$ws.Range("E51").Value = "  -2.76%  "
